$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LTPI 1.09")
$ws.Name = "LTPI 1.09 & 1.1"
